$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wXmlns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function New-PkgXml($bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wXmlns + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) Row 14 ("purchaseBeverage1"), Actual Results cell: drop the stray
#    _GoBack bookmark that trails "Coffee purchased successfully."
# ---------------------------------------------------------------------
$cell14 = $t.Cell(14, 4)
$p14 = $cell14.Range.Paragraphs.Item(1)
$r14 = $p14.Range

$body14 = '<w:p w14:paraId="01B07DC3" w14:textId="6BA67489" w:rsidR="00BC14B9" w:rsidRPr="007E191E" w:rsidRDefault="000E6ADA" w:rsidP="00BC14B9"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Coffee purchased successfully.</w:t></w:r></w:p>'
$r14.InsertXML((New-PkgXml $body14))

# ---------------------------------------------------------------------
# 2) Row 15 ("purchaseBeverage2"), Actual Results cell: it was left
#    blank - add "Coffee could not be purchased."
# ---------------------------------------------------------------------
$cell15 = $t.Cell(15, 4)
$p15 = $cell15.Range.Paragraphs.Item(1)
$r15 = $p15.Range
$r15.Collapse(1)

$body15 = '<w:p w14:paraId="41674478" w14:textId="77777777" w:rsidR="00BC14B9" w:rsidRPr="007E191E" w:rsidRDefault="00BC14B9" w:rsidP="00BC14B9"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Coffee could not be purchased.</w:t></w:r></w:p>'
$r15.InsertXML((New-PkgXml $body15))

# ---------------------------------------------------------------------
# 3) Row 16 ("purchaseBeverage3"), Actual Results cell: it was left
#    blank - add "Coffee could not be purchased." and the _GoBack
#    bookmark now lands here (last edit location).
# ---------------------------------------------------------------------
$cell16 = $t.Cell(16, 4)
$p16 = $cell16.Range.Paragraphs.Item(1)
$r16 = $p16.Range
$r16.Collapse(1)

$body16 = '<w:p w14:paraId="570FE503" w14:textId="77777777" w:rsidR="00BC14B9" w:rsidRPr="007E191E" w:rsidRDefault="00BC14B9" w:rsidP="00BC14B9"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Coffee could not be purchased.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$r16.InsertXML((New-PkgXml $body16))
